# "Agora excel tem butões"
# Renames the technician names in the "Técnicos" sheet to the generic
# "Técnico N" placeholders, renames the first two project names in
# "Projetos" to "Projeto 0"/"Projeto 1", updates both sheets' selections,
# and switches the active tab from "Técnicos" to "Projetos".

$wb = $excel.ActiveWorkbook

$wsTecnicos = $wb.Worksheets.Item("Técnicos")
$wsProjetos = $wb.Worksheets.Item("Projetos")

# --- Técnicos sheet: replace the B2:B14 technician names --------------
$tecnicoNomes = @(
    "Técnico 0",
    "Técnico 1",
    "Técnico 2",
    "Técnico 3",
    "Técnico 4",
    "Técnico 5",
    "Técnico 6",
    "Técnico 7",
    "Técnico 8",
    "Técnico 9",
    "Técnico 10",
    "Técnico 11",
    "Técnico 12"
)

for ($i = 0; $i -lt $tecnicoNomes.Length; $i++) {
    $row = 2 + $i
    $wsTecnicos.Cells.Item($row, 2).Value = $tecnicoNomes[$i]
}

# --- Projetos sheet: replace the first two (placeholder) project names -
$wsProjetos.Cells.Item(2, 2).Value = "Projeto 0"
$wsProjetos.Cells.Item(3, 2).Value = "Projeto 1"

# --- Update selections on each sheet -----------------------------------
$wsTecnicos.Range("B2:B14").Select()
$wsProjetos.Range("B2:B19").Select()

# --- Switch the active sheet to "Projetos" ------------------------------
$wsProjetos.Activate()
